# Apply Phase 3 documentation updates to the Gantt chart worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt Chart")

# --- Simple % complete updates (column H) ---
$ws.Range("H43").Value = 0.8
$ws.Range("H44").Value = 0.8
$ws.Range("H45").Value = 0.8
$ws.Range("H46").Value = 1.0
$ws.Range("H47").Value = 0.8
$ws.Range("H50").Value = 1.0
$ws.Range("H54").Value = 1.0
$ws.Range("H55").Value = 1.0
$ws.Range("H56").Value = 1.0
$ws.Range("H57").Value = 1.0
$ws.Range("H58").Value = 1.0
$ws.Range("H59").Value = 1.0
$ws.Range("H60").Value = 0.8

# --- Row 48: Implement Hand Class - Testing (Brendom) ---
$ws.Range("E48").Value = 45598.0
$ws.Range("F48").Value = 45610.0
$ws.Range("H48").Value = 0.8

# --- Row 49: Send/receive message (Carolina) ---
$ws.Range("E49").Value = 45599.0
$ws.Range("F49").Value = 45611.0
$ws.Range("H49").Value = 0.8

# --- Row 51: GUI / Login (Brendom) ---
$ws.Range("E51").Value = 45622.0
$ws.Range("F51").Value = 45632.0
$ws.Range("H51").Value = 1.0

# --- Row 55: Update SRS, owner Matthew ---
$ws.Range("D55").Value = "Matthew"
$ws.Range("E55").Value = 45623.0
$ws.Range("F55").Value = 45630.0

# --- Row 56: Update Design, dates added ---
$ws.Range("E56").Value = 45624.0
$ws.Range("F56").Value = 45631.0

# --- Row 57: Update UML Diagram, owner changed Everyone -> Matthew ---
$ws.Range("D57").Value = "Matthew"
$ws.Range("E57").Value = 45625.0
$ws.Range("F57").Value = 45632.0

# --- Row 59: Add Phase 3 to Github with all changes, owner Matthew ---
$ws.Range("D59").Value = "Matthew"
$ws.Range("E59").Value = 45623.0
$ws.Range("F59").Value = 45629.0

# --- Row 60: Merge / Testing (server/client/gui), owner Everyone ---
$ws.Range("D60").Value = "Everyone"

# --- Row 61: Demo/presentation, owner Everyone ---
$ws.Range("D61").Value = "Everyone"

# Un-highlight the week-9 M/T/W/R cells for row 61 (copy format from AH61, which
# already carries the plain style used elsewhere in that row).
$ws.Range("AH61").Copy()
$ws.Range("BC61:BF61").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 62: clear WBS number value ---
$ws.Range("B62").Value = ""

$wb.Save()
